$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "hi"

$ws.Range("A2").Value = "my"
$ws.Range("B2").Value = "name"
$ws.Range("C2").Value = "is "
$ws.Range("D2").Value = "sheela"
